# Adds the new log entry row (row 25) describing the checks functionality work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the per-column formatting used by the rows just above (date /
# plain text / wrapped text) onto the new row before filling in values.
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C25").PasteSpecial(-4122)

# New row: date, time spent, description
$ws.Range("A25").Value = 45250
$ws.Range("B25").Value = "~5 hrs"
$ws.Range("C25").Value = "added the menu and the controller that handles most of the program logic"

$ws.Rows.Item(25).RowHeight = 45

# Update the selection to reflect the new active cell after data entry
$ws.Range("C25").Select()
